$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-25 (HOUR 1-24)
# Columns: B = TOTAL_SUBSTATION_LOAD, C = CONTESTABLE_ENERGY, D = ACTUAL_ENERGY
$data = @(
    @(36728, 5799.800999999999, 30928.199),
    @(35269, 5649.136500000001, 29619.8635),
    @(33348, 5554.8325, 27793.1675),
    @(31700, 5481.98, 26218.02),
    @(30762, 5456.8815, 25305.1185),
    @(30953, 5512.395, 25440.605),
    @(31330, 5765.477227722772, 25564.52277227723),
    @(32595, 6519.4465, 26075.5535),
    @(36208, 7901.5895, 28306.4105),
    @(39069, 13052.0005, 26016.9995),
    @(39884, 14886.7985, 24997.2015),
    @(39218, 14810.005, 24407.995),
    @(38683, 14699.3385, 23983.6615),
    @(40878, 15357.636, 25520.364),
    @(41508, 15432.207, 26075.793),
    @(40766, 15338.9075, 25427.0925),
    @(38994, 15792.80217625723, 23201.19782374277),
    @(39323, 15557.50171551809, 23765.49828448191),
    @(41039, 15152.49877462994, 25886.50122537006),
    @(38646, 13812.90902852661, 24833.09097147339),
    @(38253, 12007.84432898735, 26245.15567101265),
    @(38679, 9523.143, 29155.857),
    @(38098, 6375.7855, 31722.2145)
)

$row = 2
foreach ($vals in $data) {
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $row++
}

# Row 25 only has C and D updated, B25 remains empty
$ws.Cells.Item(25, 3).Value = 5494.996500000001
$ws.Cells.Item(25, 4).Value = 0
